$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the last existing data row (A374) down through
# the new rows (A375:A385) so the new date cells pick up the same style
# (date number format / border / font) as the rest of column A.
$ws.Range("A374").Copy()
$ws.Range("A375:A385").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New data rows to append (row, date serial, col B, col C, col D)
$data = @(
    @(375, 44449, 1, 5, 152.6717557251908),
    @(376, 44450, 1, 4, 122.1374045801527),
    @(377, 44451, 0, 4, 122.1374045801527),
    @(378, 44452, 1, 4, 122.1374045801527),
    @(379, 44453, 0, 4, 122.1374045801527),
    @(380, 44454, 0, 3, 91.6030534351145),
    @(381, 44455, 0, 3, 91.6030534351145),
    @(382, 44456, 0, 2, 61.06870229007634),
    @(383, 44457, 2, 3, 91.6030534351145),
    @(384, 44458, 2, 5, 152.6717557251908),
    @(385, 44459, 0, 4, 122.1374045801527)
)

foreach ($row in $data) {
    $r = $row[0]
    $dateSerial = $row[1]
    $b = $row[2]
    $c = $row[3]
    $d = $row[4]

    $ws.Cells.Item($r, 1).Value = $dateSerial
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
}
